$wb = $excel.ActiveWorkbook
$neg = $wb.Worksheets.Item("NegativeTests")
$neg.Range("B7").Value2 = "test@test.com"
$neg.Hyperlinks.Add($neg.Range("B7:B7"), "mailto:test@test.com") | Out-Null
